# Updated team report to reflect next sprint
$wb = $excel.ActiveWorkbook

# --- Backlog: fill in the assigned owner (initials) for the Sprint 3
#     stories (US13-US18), matching the pattern already used for the
#     stories of earlier sprints.
$backlog = $wb.Worksheets.Item("Backlog")
$backlog.Range("D14").Value = "mm"
$backlog.Range("D15").Value = "mm"
$backlog.Range("D16").Value = "bg"
$backlog.Range("D17").Value = "bg"
$backlog.Range("D18").Value = "rh"
$backlog.Range("D19").Value = "rh"

# --- Sprint3: populate the next sprint's worksheet with the stories
#     pulled in from the backlog, their owners, status and estimates.
$sprint3 = $wb.Worksheets.Item("Sprint3")

$sprint3.Range("A2").Value = "US13"
$sprint3.Range("B2").Value = "Siblings spacing"
$sprint3.Range("C2").Value = "mm"
$sprint3.Range("D2").Value = "assigned"
$sprint3.Range("E2").Value = 30
$sprint3.Range("F2").Value = 20

$sprint3.Range("A3").Value = "US14"
$sprint3.Range("B3").Value = "Multiple births less than 5"
$sprint3.Range("C3").Value = "mm"
$sprint3.Range("D3").Value = "assigned"
$sprint3.Range("E3").Value = 30
$sprint3.Range("F3").Value = 20

$sprint3.Range("A4").Value = "US15"
$sprint3.Range("B4").Value = "Fewer than 15 siblings"
$sprint3.Range("C4").Value = "bg"
$sprint3.Range("D4").Value = "assigned"
$sprint3.Range("E4").Value = 10
$sprint3.Range("F4").Value = 15

$sprint3.Range("A5").Value = "US16"
$sprint3.Range("B5").Value = "Male last names"
$sprint3.Range("C5").Value = "bg"
$sprint3.Range("D5").Value = "assigned"
$sprint3.Range("E5").Value = 45
$sprint3.Range("F5").Value = 45

$sprint3.Range("A6").Value = "US17"
$sprint3.Range("B6").Value = "No marriages to descendants"
$sprint3.Range("C6").Value = "rh"
$sprint3.Range("D6").Value = "assigned"
$sprint3.Range("E6").Value = 45
$sprint3.Range("F6").Value = 45

$sprint3.Range("A7").Value = "US18"
$sprint3.Range("B7").Value = "Siblings should not marry"
$sprint3.Range("C7").Value = "rh"
$sprint3.Range("D7").Value = "assigned"
$sprint3.Range("E7").Value = 30
$sprint3.Range("F7").Value = 20

# --- Restore the selections the author left on each sheet.
$backlog.Range("B14").Select()
$wb.Worksheets.Item("Burndown").Range("F4").Select()
$wb.Worksheets.Item("Sprint2").Range("I5").Select()
$sprint3.Range("J1").Select()
$wb.Worksheets.Item("Stories").Range("B15").Select()

$backlog.Activate()
$sprint3.Activate()
$wb.Worksheets.Item("Sprint2").Activate()
